# Add data for 2022-03-01
# - Renames the sheet / report-date label from "...February 20" to "...February 21"
# - Updates carjacking counts for several neighborhoods/months to reflect newly
#   reconciled data (including a few brand-new, previously-empty data points).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab to reflect the new "through" date
$ws.Name = "Through 2022-02-21"

# Update the column header label (shared string) to match the new "through" date
$ws.Range("B1").Value = "February 2022 (through February 21)"

# Updated / newly-added neighborhood counts
$ws.Range("F2").Value = 3     # Englewood, February 2020
$ws.Range("D3").Value = 10    # Austin, February 2021
$ws.Range("F3").Value = 7     # Austin, February 2020
$ws.Range("B4").Value = 2     # New City, February 2022 (through Feb 21)
$ws.Range("B6").Value = 5     # South Shore, February 2022 (through Feb 21)
$ws.Range("D6").Value = 8     # South Shore, February 2021
$ws.Range("F8").Value = 6     # North Lawndale, February 2020
$ws.Range("D11").Value = 8    # Garfield Park, February 2021
$ws.Range("H11").Value = 2    # Garfield Park, February 2019
$ws.Range("P11").Value = 1    # Garfield Park, February 2015 (new value)
$ws.Range("D13").Value = 6    # Little Italy, UIC, February 2021
$ws.Range("B17").Value = 3    # Chatham, February 2022 (through Feb 21)
$ws.Range("D17").Value = 2    # Chatham, February 2021
$ws.Range("B20").Value = 2    # Kenwood, February 2022 (through Feb 21)
$ws.Range("B22").Value = 2    # Humboldt Park, February 2022 (through Feb 21)
$ws.Range("J23").Value = 1    # Washington Heights, February 2018 (new value)
$ws.Range("F34").Value = 2    # West Loop, February 2020
$ws.Range("D36").Value = 2    # South Chicago, February 2021
$ws.Range("B42").Value = 2    # River North, February 2022 (through Feb 21)
$ws.Range("B46").Value = 2    # Lower West Side, February 2022 (through Feb 21)
$ws.Range("D47").Value = 3    # Little Village, February 2021
$ws.Range("B48").Value = 4    # Lake View, February 2022 (through Feb 21)
$ws.Range("J85").Value = 1    # Washington Park, February 2018 (new value)
$ws.Range("L85").Value = 2    # Washington Park, February 2017
